$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'heba@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Veronia.rafat@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G3").Value = 'backup@backdoor.com, gehanadel@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg'
$ws.Range("G4").Value = 'heba@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G5").Value = 'AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range("G17").Value = 'user@user.com, enas.omran@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Range("G18").Value = 'eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G20").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range("G21").Value = 'heba@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Veronia.rafat@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G22").Value = 'backup@backdoor.com, gehanadel@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg'
$ws.Range("G23").Value = 'heba@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G24").Value = 'Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg'
$ws.Range("G28").Value = 'arwaelsayed03@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, Sarah.Abdelmohsen@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg'
$ws.Range("G29").Value = 'yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G36").Value = 'user@user.com, enas.omran@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Range("G37").Value = 'eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G38").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range("G39").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range("G40").Value = 'heba@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Veronia.rafat@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G41").Value = 'hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G42").Value = 'shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G43").Value = 'Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg'
$ws.Range("G47").Value = 'arwaelsayed03@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, merna.said@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G48").Value = 'Sarah.Abdelmohsen@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G54").Value = 'afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg'
$ws.Range("G56").Value = 'eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G57").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range("G58").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range("G59").Value = 'heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G60").Value = 'hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G63").Value = 'aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G66").Value = 'dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G67").Value = 'yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G75").Value = 'eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G77").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range("G78").Value = 'heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G79").Value = 'hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G81").Value = 'Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg'
$ws.Range("G83").Value = 'aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G85").Value = 'dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G86").Value = 'Sarah.Abdelmohsen@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg'
$ws.Range("G94").Value = 'eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G96").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range("G97").Value = 'heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G98").Value = 'hend_mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, hananragab@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg'
$ws.Range("G99").Value = 'shaimaa.ahmed@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G100").Value = 'Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg'
$ws.Range("G101").Value = 'aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G104").Value = 'arwaelsayed03@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, merna.said@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg'
$ws.Range("G113").Value = 'eman.samir@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range("G115").Value = 'marinasorial@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
